$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 46328.5335549085
$ws.Range("C2").Value = 96865.93318510834
$ws.Range("D2").Value = 100676.8686432669
